# "Changed GPA as number"
#
# The source workbook used the shared string "--" as a placeholder in the
# GPA column (Y) for every student who failed (column X = "F"). This
# script replaces that text placeholder with the literal number 0 so the
# GPA column is fully numeric. Once no cell references the "--" shared
# string any more, it naturally drops out of the saved shared-strings
# table (and every other shared-string reference renumbers accordingly --
# that renumbering is an automatic side effect of the save, not something
# this script needs to special-case).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$fixedCount = 0
for ($r = 2; $r -le 198; $r++) {
    $cell = $ws.Cells.Item($r, 25)   # column Y = GPA
    if ($cell.Value() -eq "--") {
        $cell.Value = 0
        $fixedCount = $fixedCount + 1
    }
}
Write-Host "GPA placeholders replaced with 0:" $fixedCount

# Restore the active selection to where the author last left off (Y4),
# scrolled so column J is leftmost in the viewport.
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 10
$ws.Range("Y4").Select()
$win.ScrollRow = 1
$win.ScrollColumn = 10
